$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "active"
$ws.Name = "active"

# Clear the (mostly cosmetic border/fill) formatting that was applied to the
# header row and the "X" marker cells - they revert to the default style.
$ws.Range("A1:D1").ClearFormats()
$ws.Range("B2").ClearFormats()
$ws.Range("C2").ClearFormats()
$ws.Range("B3").ClearFormats()
$ws.Range("D3").ClearFormats()

# The two empty, formatted-only cells disappear entirely.
$ws.Range("D2").Clear()
$ws.Range("C3").Clear()

# The date cells in column A keep a date number format, but lose the
# border formatting that used to be bundled with it.
$ws.Range("A2:A4").ClearFormats()
$ws.Range("A2:A4").NumberFormat = "m/d/yy"

# Move the active selection.
$ws.Range("D21").Select()
